$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert two new rows at the top of the data (rows 2 and 3), pushing the
# existing 18 IPO rows down to rows 4-21.
$ws1.Rows.Item(2).Insert()
$ws1.Rows.Item(2).Insert()

# Excel auto-detects "YYYY-MM-DD"-shaped text as a date; force columns A-C
# (date columns) on the two new rows to stay plain text, matching every
# other row in the sheet (which stores these as shared-string text, not
# real dates).
$ws1.Range("A2:C3").NumberFormat = "@"

# Row 2: NH SPAC (엔에이치스팩31호)
$ws1.Cells.Item(2,1).Value = '2024-07-09'
$ws1.Cells.Item(2,2).Value = '2024-07-10'
$ws1.Cells.Item(2,3).Value = '2024-07-26'
$ws1.Cells.Item(2,4).Value = 'NH'
$ws1.Cells.Item(2,5).Value = '엔에이치스팩31호'
$ws1.Cells.Item(2,6).Value = 6000000
$ws1.Cells.Item(2,7).Value = 6000000
$ws1.Cells.Item(2,8).Value = 0
$ws1.Cells.Item(2,9).Value = 2000
$ws1.Cells.Item(2,10).Value = 2000
$ws1.Cells.Item(2,11).Value = 6345000
$ws1.Cells.Item(2,12).Value = 0
$ws1.Cells.Item(2,13).Value = 2000
$ws1.Cells.Item(2,14).Value = '1123.43 :1'
$ws1.Cells.Item(2,15).Value = '-'
$ws1.Cells.Item(2,16).Value = 0
$ws1.Cells.Item(2,17).Value = 0
$ws1.Cells.Item(2,18).Value = 0
$ws1.Cells.Item(2,19).Value = 0
$ws1.Cells.Item(2,20).Value = 0
$ws1.Cells.Item(2,21).Value = 0
$ws1.Cells.Item(2,22).Value = 0
$ws1.Cells.Item(2,23).Value = 0
$ws1.Cells.Item(2,24).Value = 0
$ws1.Cells.Item(2,25).Value = '금융지원서비스업'

# Row 3: SK SPAC (SK증권제13호스팩)
$ws1.Cells.Item(3,1).Value = '2024-07-09'
$ws1.Cells.Item(3,2).Value = '2024-07-10'
$ws1.Cells.Item(3,3).Value = '2024-07-25'
$ws1.Cells.Item(3,4).Value = 'SK'
$ws1.Cells.Item(3,5).Value = 'SK증권제13호스팩'
$ws1.Cells.Item(3,6).Value = 4000000
$ws1.Cells.Item(3,7).Value = 4000000
$ws1.Cells.Item(3,8).Value = 0
$ws1.Cells.Item(3,9).Value = 2000
$ws1.Cells.Item(3,10).Value = 2000
$ws1.Cells.Item(3,11).Value = 4420000
$ws1.Cells.Item(3,12).Value = 0
$ws1.Cells.Item(3,13).Value = 2000
$ws1.Cells.Item(3,14).Value = '1197.45:1'
$ws1.Cells.Item(3,15).Value = '-'
$ws1.Cells.Item(3,16).Value = 0
$ws1.Cells.Item(3,17).Value = 0
$ws1.Cells.Item(3,18).Value = 0
$ws1.Cells.Item(3,19).Value = 0
$ws1.Cells.Item(3,20).Value = 0
$ws1.Cells.Item(3,21).Value = 0
$ws1.Cells.Item(3,22).Value = 0
$ws1.Cells.Item(3,23).Value = 0
$ws1.Cells.Item(3,24).Value = 0
$ws1.Cells.Item(3,25).Value = '기업인수목적 주식회사'

# Drop the formatting Insert()/NumberFormat picked up so the two new rows
# end up with the sheet's default (unstyled) cells, same as every other
# data row.
$ws1.Range("A2:Y3").ClearFormats()
